$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.06476221978664398
$ws.Cells.Item(2, 2).Value = 0.9810612201690674
$ws.Cells.Item(2, 3).Value = 0.00437135249376297
$ws.Cells.Item(2, 4).Value = 0.9992126822471619
$ws.Cells.Item(3, 1).Value = 0.009770811535418034
$ws.Cells.Item(3, 2).Value = 0.9983066320419312
$ws.Cells.Item(3, 3).Value = 0.002968936692923307
$ws.Cells.Item(3, 4).Value = 0.9993557929992676
$ws.Cells.Item(4, 1).Value = 0.004503690171986818
$ws.Cells.Item(4, 2).Value = 0.9989683628082275
$ws.Cells.Item(4, 3).Value = 0.0005921124829910696
$ws.Cells.Item(4, 4).Value = 0.9996421337127686
$ws.Cells.Item(5, 1).Value = 0.002431630855426192
$ws.Cells.Item(5, 2).Value = 0.9993966221809387
$ws.Cells.Item(5, 3).Value = 0.0007437243475578725
$ws.Cells.Item(5, 4).Value = 0.9998568296432495
$ws.Cells.Item(6, 1).Value = 0.001630124868825078
$ws.Cells.Item(6, 2).Value = 0.9996301531791687
$ws.Cells.Item(6, 3).Value = 0.0007195430225692689
$ws.Cells.Item(6, 4).Value = 0.9998568296432495
$ws.Cells.Item(7, 1).Value = 0.001827004365622997
$ws.Cells.Item(7, 2).Value = 0.9997080564498901
$ws.Cells.Item(7, 3).Value = 0.0003458843566477299
$ws.Cells.Item(7, 4).Value = 0.9997852444648743
$ws.Cells.Item(8, 1).Value = 0.001015970483422279
$ws.Cells.Item(8, 2).Value = 0.9997080564498901
$ws.Cells.Item(8, 3).Value = 0.0002899055252783
$ws.Cells.Item(8, 4).Value = 0.9998568296432495
$ws.Cells.Item(9, 1).Value = 0.001161716412752867
$ws.Cells.Item(9, 2).Value = 0.9996885657310486
$ws.Cells.Item(9, 3).Value = 0.00005530666749109514
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(10, 1).Value = 0.0004956370103172958
$ws.Cells.Item(10, 2).Value = 0.9998443126678467
$ws.Cells.Item(10, 3).Value = 0.0001134614940383472
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(11, 1).Value = 0.0008988276240415871
$ws.Cells.Item(11, 2).Value = 0.9998053312301636
$ws.Cells.Item(11, 3).Value = 0.00005356170368031599
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 1).Value = 0.0008424934349022806
$ws.Cells.Item(12, 2).Value = 0.9998248219490051
$ws.Cells.Item(12, 3).Value = 0.00001095504740078468
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 1).Value = 0.0001231416244991124
$ws.Cells.Item(13, 2).Value = 0.9999610781669617
$ws.Cells.Item(13, 3).Value = 0.000001390909119436401
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(14, 1).Value = 0.0008421097882091999
$ws.Cells.Item(14, 2).Value = 0.9998248219490051
$ws.Cells.Item(14, 3).Value = 0.0001134365884354338
$ws.Cells.Item(14, 4).Value = 0.9999284148216248
$ws.Cells.Item(15, 1).Value = 0.0004102381353732198
$ws.Cells.Item(15, 2).Value = 0.9999221563339233
$ws.Cells.Item(15, 3).Value = 0.000001188896476378432
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(16, 1).Value = 0.0003576846211217344
$ws.Cells.Item(16, 2).Value = 0.9999026656150818
$ws.Cells.Item(16, 3).Value = 0.00001221865022671409
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(17, 1).Value = 0.0001032729560392909
$ws.Cells.Item(17, 2).Value = 0.9999415874481201
$ws.Cells.Item(17, 3).Value = 0.0000006898830520185584
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(18, 1).Value = 0.0006652430165559053
$ws.Cells.Item(18, 2).Value = 0.9998443126678467
$ws.Cells.Item(18, 3).Value = 0.0000002514076413717703
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 1).Value = 0.0002929032489191741
$ws.Cells.Item(19, 2).Value = 0.9999026656150818
$ws.Cells.Item(19, 3).Value = 0.000001966115860341233
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = 0.0004951423616148531
$ws.Cells.Item(20, 2).Value = 0.9998443126678467
$ws.Cells.Item(20, 3).Value = 0.0000008800275850262551
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(21, 1).Value = 0.0001083862953237258
$ws.Cells.Item(21, 2).Value = 0.9999610781669617
$ws.Cells.Item(21, 3).Value = 0.00000002609093918692906
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(22, 1).Value = 0.0002073797077173367
$ws.Cells.Item(22, 2).Value = 0.9999415874481201
$ws.Cells.Item(22, 3).Value = 0.0000003957772207741073
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(23, 1).Value = 0.0006078690057620406
$ws.Cells.Item(23, 2).Value = 0.9998248219490051
$ws.Cells.Item(23, 3).Value = 0.0009092004620470107
$ws.Cells.Item(23, 4).Value = 0.9998568296432495
$ws.Cells.Item(24, 1).Value = 0.0006343477289192379
$ws.Cells.Item(24, 2).Value = 0.9998248219490051
$ws.Cells.Item(24, 3).Value = 0.00003029372965102084
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(25, 1).Value = 0.0001211665221489966
$ws.Cells.Item(25, 2).Value = 0.9999610781669617
$ws.Cells.Item(25, 3).Value = 0.000002032799329754198
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 1).Value = 0.00003348788595758379
$ws.Cells.Item(26, 3).Value = 0.0000001612332312106446
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.00004301036096876487
$ws.Cells.Item(27, 2).Value = 0.9999805092811584
$ws.Cells.Item(27, 3).Value = 0.00001500571124779526
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = 0.0004061237850692123
$ws.Cells.Item(28, 2).Value = 0.9999221563339233
$ws.Cells.Item(28, 3).Value = 0.000000001066577048902673
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 1).Value = 0.0005263984203338623
$ws.Cells.Item(29, 2).Value = 0.999883234500885
$ws.Cells.Item(29, 3).Value = 0.0000003645078834324522
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = 0.0003594472364056855
$ws.Cells.Item(30, 2).Value = 0.9999415874481201
$ws.Cells.Item(30, 3).Value = 0.00000001127131277200988
$ws.Cells.Item(30, 4).Value = 1
$ws.Cells.Item(31, 1).Value = 0.00008497395901940763
$ws.Cells.Item(31, 2).Value = 0.9999805092811584
$ws.Cells.Item(31, 3).Value = 0.000000001783313718917157
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.000004485066710913088
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 0.000000001066576271746555
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.000006852985734440153
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 3).Value = 0.000000000759402374317375
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = 0.0000006257054678826535
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = 0.0000000003413047622302656
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 1).Value = 0.0001588651939528063
$ws.Cells.Item(35, 2).Value = 0.9999221563339233
$ws.Cells.Item(35, 3).Value = 0.00000009158576830259335
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(36, 1).Value = 0.001103704678826034
$ws.Cells.Item(36, 2).Value = 0.9998248219490051
$ws.Cells.Item(36, 3).Value = 0.0000000312516306166799
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(37, 1).Value = 0.0001368596276734024
$ws.Cells.Item(37, 2).Value = 0.9999610781669617
$ws.Cells.Item(37, 3).Value = 0.0000000001706524505040719
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = 0.0002708659158088267
$ws.Cells.Item(38, 2).Value = 0.9999415874481201
$ws.Cells.Item(38, 3).Value = 0.0000000008703239240404059
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = 0.00009745963325258344
$ws.Cells.Item(39, 2).Value = 0.9999805092811584
$ws.Cells.Item(39, 3).Value = 0.000000000537554445401156
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 1).Value = 0.0000102037347460282
$ws.Cells.Item(40, 2).Value = 1
$ws.Cells.Item(40, 3).Value = 0.00000000006826097326273484
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 1).Value = 0.0001463483495172113
$ws.Cells.Item(41, 2).Value = 0.9999415874481201
$ws.Cells.Item(41, 3).Value = 0.000000001168964369746561
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(42, 1).Value = 0.0007930033607408404
$ws.Cells.Item(42, 2).Value = 0.9998443126678467
$ws.Cells.Item(42, 3).Value = 0.0000000003071743415183903
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 1).Value = 0.0002277269959449768
$ws.Cells.Item(43, 2).Value = 0.9999221563339233
$ws.Cells.Item(43, 3).Value = 0.0000000001450545655812263
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(44, 1).Value = 0.00003462206950644031
$ws.Cells.Item(44, 2).Value = 0.9999805092811584
$ws.Cells.Item(44, 3).Value = 0.0000000005034239691781295
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(45, 1).Value = 0.00009538239100947976
$ws.Cells.Item(45, 2).Value = 0.9999805092811584
$ws.Cells.Item(45, 3).Value = 0.00000003741920906463747
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 1).Value = 0.0003430102369748056
$ws.Cells.Item(46, 2).Value = 0.9999415874481201
$ws.Cells.Item(46, 3).Value = 0.0000000002901090201401502
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(47, 1).Value = 0.000007665503289899789
$ws.Cells.Item(47, 2).Value = 1
$ws.Cells.Item(47, 3).Value = 0.00003250773079344071
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(48, 1).Value = 0.0005710429977625608
$ws.Cells.Item(48, 2).Value = 0.9999026656150818
$ws.Cells.Item(48, 3).Value = 0.0004814779968000948
$ws.Cells.Item(48, 4).Value = 0.9999284148216248
$ws.Cells.Item(49, 1).Value = 0.0002022783446591347
$ws.Cells.Item(49, 2).Value = 0.9999805092811584
$ws.Cells.Item(49, 3).Value = 0.00000003297953199421499
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = 0.0001402939815307036
$ws.Cells.Item(50, 2).Value = 0.9999805092811584
$ws.Cells.Item(50, 3).Value = 0.0000001184031077627878
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 1).Value = 0.0002179859584430233
$ws.Cells.Item(51, 2).Value = 0.9999221563339233
$ws.Cells.Item(51, 3).Value = 0.0000000001791850418042529
$ws.Cells.Item(51, 4).Value = 1
